$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: replace the helpline resource row with the AA Meetings / AOD resource row ---
$ws.Range("A2").Value = "AOD_Aameetings"
$ws.Range("B2").Value = "AA Meetings"
$ws.Range("C2").Value = "AODdata"
$ws.Range("D2").Value = "AA Meetings"
$ws.Range("E2").Value = "F"
$ws.Range("F2").Value = "T"

# B2 and D2 switch from the "text" style (s=1) to the "general" style used by A2 (s=3)
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Clear the old standalone labels that used to live in A11:A13 ---
$ws.Range("A11").ClearContents()
$ws.Range("A12").ClearContents()
$ws.Range("A13").Clear()

# --- Reduce the trailing blank block: keep C17 (style 1) and C18 (style 4, was C22's style) ---
$ws.Range("C22").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C19:C22").Clear()

# --- New AOD resource "ui" lines (rows 20-23) ---
$ws.Range("A20").Value = "mod_Accordion_ui('AOD_NaloxoneCal')"
$ws.Range("A21").Value = "mod_Accordion_ui('AOD_SMART')"
$ws.Range("A22").Value = "mod_Accordion_ui('AOD_Refuge')"
$ws.Range("A23").Value = "mod_Accordion_ui('AOD_Aameetings')"

# --- New AOD resource "server"/"info" line pairs (rows 26-33) ---
$ws.Range("A26").Value = "mod_Accordion_server('AOD_NaloxoneCal', selector=selection, data=AODdata, title = c('OASAS Naloxone Training Calendar'), Visible = T)"
$ws.Range("A27").Value = "mod_info_server('AOD_NaloxoneCal', selector = selection, data = AODdata, rownametitle = c('OASAS Naloxone Training Calendar'), phone = F, website = T)"
$ws.Range("A28").Value = "mod_Accordion_server('AOD_SMART', selector=selection, data=AODdata, title = c('SMART Recovery Meetings'), Visible = T)"
$ws.Range("A29").Value = "mod_info_server('AOD_SMART', selector = selection, data = AODdata, rownametitle = c('SMART Recovery Meetings'), phone = F, website = T)"
$ws.Range("A30").Value = "mod_Accordion_server('AOD_Refuge', selector=selection, data=AODdata, title = c('Refuge Recovery Meetings'), Visible = T)"
$ws.Range("A31").Value = "mod_info_server('AOD_Refuge', selector = selection, data = AODdata, rownametitle = c('Refuge Recovery Meetings'), phone = F, website = T)"
$ws.Range("A32").Value = "mod_Accordion_server('AOD_Aameetings', selector=selection, data=AODdata, title = c('AA Meetings'), Visible = T)"
$ws.Range("A33").Value = "mod_info_server('AOD_Aameetings', selector = selection, data = AODdata, rownametitle = c('AA Meetings'), phone = F, website = T)"

# New cells use the same "general" style (s=3) as A2/A13
$ws.Range("A2").Copy()
$ws.Range("A20:A23").PasteSpecial(-4122)
$ws.Range("A26:A33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the workbook/sheet view: scrolled to A13, selection A17:XFD20 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("A17:XFD20").Select()

$wb.Save()
